$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string (rich text) partial edits ---
$ws.Range("A8").Characters(21, 2).Text = "19"
$ws.Range("C9").Characters(27, 9).Text = "5/5/2025"
$ws.Range("C9").Characters(46, 8).Text = "5/11/2025"

# --- Column E width (closest achievable to 7.433768 given COM quantization) ---
$ws.Columns("E").ColumnWidth = 6.71

# --- Cell value / style updates ---
$ws.Range("L14").Value = 0
$ws.Range("L14").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = -42.857142857142
$ws.Range("N15").Value = -73.333333333333
$ws.Range("D16").Value = 3
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = -33.333333333333
$ws.Range("E16").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 125
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = 88.888888888888
$ws.Range("L16").Value = 3.030303030303
$ws.Range("M16").Value = -68.518518518518
$ws.Range("N16").Value = -90.229885057471
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 75
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = -25.742574257425
$ws.Range("L17").Value = -8.536585365853
$ws.Range("M17").Value = -26.470588235294
$ws.Range("N17").Value = -34.210526315789
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = 21.27659574468
$ws.Range("L18").Value = 1.785714285714
$ws.Range("M18").Value = -55.11811023622
$ws.Range("N18").Value = -90.25641025641
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 34.615384615384
$ws.Range("I19").Value = 141
$ws.Range("J19").Value = 144
$ws.Range("K19").Value = -2.083333333333
$ws.Range("L19").Value = 3.676470588235
$ws.Range("M19").Value = -14.545454545454
$ws.Range("N19").Value = -28.061224489795
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = -54.545454545454
$ws.Range("I20").Value = 63
$ws.Range("J20").Value = 81
$ws.Range("K20").Value = -22.222222222222
$ws.Range("L20").Value = 10.526315789473
$ws.Range("M20").Value = -51.908396946564
$ws.Range("N20").Value = -94.810543657331
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 78
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = -4.878048780487
$ws.Range("I21").Value = 375
$ws.Range("J21").Value = 397
$ws.Range("K21").Value = -5.541561712846
$ws.Range("L21").Value = 1.902173913043
$ws.Range("M21").Value = -41.950464396284
$ws.Range("N21").Value = -84.866828087167
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = -15.384615384615
$ws.Range("I24").Value = 192
$ws.Range("J24").Value = 249
$ws.Range("K24").Value = -22.89156626506
$ws.Range("L24").Value = -31.182795698924
$ws.Range("M24").Value = -37.662337662337
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 8.333333333333
$ws.Range("I25").Value = 46
$ws.Range("J25").Value = 56
$ws.Range("K25").Value = -17.857142857142
$ws.Range("L25").Value = -9.803921568627
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -57.142857142857
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -3.125
$ws.Range("I26").Value = 155
$ws.Range("J26").Value = 124
$ws.Range("K26").Value = 25
$ws.Range("L26").Value = 31.355932203389
$ws.Range("M26").Value = -31.111111111111
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = -55.555555555555
$ws.Range("L27").Value = -20
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 100
$ws.Range("E28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 17
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 6.25
$ws.Range("L28").Value = 183.333333333333
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("G31").Value = 1
$ws.Range("G31").NumberFormat = "#,##0"
$ws.Range("H31").Value = 0
$ws.Range("H31").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = 0
